# Weekly refresh: insert a new 6-row batch (Murcott, date 44826) at the top
# of the "1000" block, pushing all subsequent rows down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at row 1000 (shifts old rows 1000-1077 down to 1006-1083)
$ws.Range("A1000:T1005").EntireRow.Insert()

# Common (constant) column values for this sheet/subset
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$prodId    = 100102
$prod      = "Cítricos"
$catId     = 100102004
$cat       = "Mandarina"

# New batch rows (row, fecha, variedad, calidad, volumen, precioMin/Max/Prom, unidad, origen, precioKg, kgUnidad)
$newRows = @(
    @{ Row=1000; Fecha=44826; Variedad="Murcott"; Calidad="Especial"; Volumen=75; Precio=7000;  Unidad="`$/bandeja 10 kilos"; Origen="Provincia de Quillota"; PrecioKg=700; KgUnidad=10 },
    @{ Row=1001; Fecha=44826; Variedad="Murcott"; Calidad="Especial"; Volumen=70; Precio=10000; Unidad="`$/caja 15 kilos";     Origen="Provincia de Quillota"; PrecioKg=667; KgUnidad=15 },
    @{ Row=1002; Fecha=44826; Variedad="Murcott"; Calidad="Primera";  Volumen=80; Precio=6000;  Unidad="`$/bandeja 10 kilos"; Origen="Provincia de Quillota"; PrecioKg=600; KgUnidad=10 },
    @{ Row=1003; Fecha=44826; Variedad="Murcott"; Calidad="Primera";  Volumen=68; Precio=9000;  Unidad="`$/caja 15 kilos";     Origen="Provincia de Quillota"; PrecioKg=600; KgUnidad=15 },
    @{ Row=1004; Fecha=44826; Variedad="Murcott"; Calidad="Segunda";  Volumen=80; Precio=5000;  Unidad="`$/bandeja 10 kilos"; Origen="Provincia de Quillota"; PrecioKg=500; KgUnidad=10 },
    @{ Row=1005; Fecha=44826; Variedad="Murcott"; Calidad="Segunda";  Volumen=65; Precio=7500;  Unidad="`$/caja 15 kilos";     Origen="Provincia de Quillota"; PrecioKg=500; KgUnidad=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value2 = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $prod
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $cat
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Precio
    $ws.Cells.Item($row, 15).Value = $r.Precio
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}

Write-Host "Done inserting new batch rows 1000-1005; dimension now through row 1083"
